# Update cryptos list values (GitHub Actions scheduled data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, avoiding Excel's automatic
# number/date inference for strings that look numeric (e.g. "308.33"),
# and make sure no stray cell style/number-format is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "46.556.04"
Set-TextValue $ws.Range("E2") "  +0.09%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.613.36"

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  +0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "308.33"
Set-TextValue $ws.Range("E5") "  +3.94%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "100.61"
Set-TextValue $ws.Range("E6") "  +2.99%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.606"
Set-TextValue $ws.Range("E7") "  +5.53%  "

# Row 8 - USDC
Set-TextValue $ws.Range("E8") "  +0.10%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.581"
Set-TextValue $ws.Range("E9") "  +13.03%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "39.56"
Set-TextValue $ws.Range("E10") "  +11.51%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("E11") "  +7.74%  "

# Row 12 - OKB
Set-TextValue $ws.Range("D12") "54.62"
Set-TextValue $ws.Range("E12") "  +1.82%  "

# Row 13 - Polkadot
Set-TextValue $ws.Range("D13") "8.20"
Set-TextValue $ws.Range("E13") "  +14.49%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "3.010.33"
Set-TextValue $ws.Range("E14") "  +7.30%  "

# Row 15 - TRON
Set-TextValue $ws.Range("E15") "  +1.61%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.618.27"
Set-TextValue $ws.Range("E16") "  +7.81%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("E17") "  +8.64%  "

# Row 18 - Chainlink
Set-TextValue $ws.Range("D18") "15.07"
Set-TextValue $ws.Range("E18") "  +7.06%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "46.714.36"
Set-TextValue $ws.Range("E19") "  +0.70%  "

# Row 20 - ShibaInu
Set-TextValue $ws.Range("E20") "  +7.35%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("E21") "  +2.39%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "6.82"
Set-TextValue $ws.Range("E22") "  +8.86%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "280.20"
Set-TextValue $ws.Range("E23") "  +13.83%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "72.08"
Set-TextValue $ws.Range("E24") "  +6.59%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "3.05"
Set-TextValue $ws.Range("E25") "  +8.33%  "

# Row 26 - ImmutableX
Set-TextValue $ws.Range("E26") "  +10.75%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "28.96"
Set-TextValue $ws.Range("E27") "  +34.75%  "

# Row 28 - Dai
Set-TextValue $ws.Range("E28") "  -0.10%  "

# Row 29 - LEO
Set-TextValue $ws.Range("D29") "4.03"
Set-TextValue $ws.Range("E29") "  +0.22%  "

# Row 30 - Cosmos
Set-TextValue $ws.Range("D30") "10.69"
Set-TextValue $ws.Range("E30") "  +8.91%  "

# Row 31 - Toncoin
Set-TextValue $ws.Range("D31") "2.31"
Set-TextValue $ws.Range("E31") "  +4.24%  "

# Row 32 - InjectiveProtocol
Set-TextValue $ws.Range("D32") "39.24"
Set-TextValue $ws.Range("E32") "  -2.33%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "6.38"
Set-TextValue $ws.Range("E33") "  +14.09%  "

# Row 34 - LidoDAOToken
Set-TextValue $ws.Range("D34") "3.67"
Set-TextValue $ws.Range("E34") "  -3.59%  "

# Rows 35/36 - WEMIXToken and Hedera swap places
Set-TextValue $ws.Range("B35") "Hedera"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D35") "0.0845"
Set-TextValue $ws.Range("E35") "  +9.09%  "

Set-TextValue $ws.Range("B36") "WEMIXToken"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D36") "2.84"
Set-TextValue $ws.Range("E36") "  +3.90%  "

# Row 37 - ARBITRUM
Set-TextValue $ws.Range("D37") "2.22"
Set-TextValue $ws.Range("E37") "  +10.57%  "

# Row 38 - Monero
Set-TextValue $ws.Range("D38") "151.04"
Set-TextValue $ws.Range("E38") "  +2.07%  "

# Row 40 - Stellar
Set-TextValue $ws.Range("E40") "  +5.83%  "

# Row 41 - EnergySwap
Set-TextValue $ws.Range("D41") "23.22"
Set-TextValue $ws.Range("E41") "  +40.63%  "

# Row 42 - Celestia
Set-TextValue $ws.Range("E42") "  +5.54%  "

# Row 43 - VeChain
Set-TextValue $ws.Range("E43") "  +10.57%  "

# Row 44 - NEARProtocol
Set-TextValue $ws.Range("E44") "  +11.16%  "

# Row 45 - RenderToken
Set-TextValue $ws.Range("E45") "  +3.96%  "

# Row 46 - Maker
Set-TextValue $ws.Range("D46") "2.144.34"
Set-TextValue $ws.Range("E46") "  +8.23%  "

# Row 47 - FirstDigitalUSD
Set-TextValue $ws.Range("D47") "0.997"
Set-TextValue $ws.Range("E47") "  -0.16%  "

# Row 48 - BitcoinSV
Set-TextValue $ws.Range("D48") "93.15"
Set-TextValue $ws.Range("E48") "  +0.31%  "

# Row 49 - FraxShare
Set-TextValue $ws.Range("E49") "  +11.69%  "

# Row 50 - Stacks
Set-TextValue $ws.Range("D50") "1.80"
Set-TextValue $ws.Range("E50") "  -2.38%  "

# Row 51 - Aave
Set-TextValue $ws.Range("D51") "109.82"
Set-TextValue $ws.Range("E51") "  +8.24%  "
